$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1822
$ws.Range("I18").Value = 799.75
$ws.Range("K18").Value = 799.75
$ws.Range("M18").Value = -515.75
$ws.Range("H69").Value = 19694.576
$ws.Range("I69").Value = 17359
$ws.Range("K69").Value = 52077
$ws.Range("M69").Value = -51203
$ws.Range("H72").Value = 19694.576
$ws.Range("I72").Value = 17359
$ws.Range("K72").Value = 156231
$ws.Range("M72").Value = -151863
$ws.Range("H87").Value = 40000
$ws.Range("J87").Value = 40000
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42496
$ws.Range("H90").Value = 40000
$ws.Range("J90").Value = 40000
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -132480
$ws.Range("H132").Value = 4042.6904
$ws.Range("I132").Value = 4579.276
$ws.Range("K132").Value = 13737.828
$ws.Range("M132").Value = -11207.828
$ws.Range("H138").Value = 4397.6294
$ws.Range("J138").Value = 5120.3335
$ws.Range("L138").Value = 15361.0005
$ws.Range("N138").Value = -25641.0005
$ws.Range("H141").Value = 5112.6665
$ws.Range("I141").Value = 4745.385
$ws.Range("K141").Value = 14236.155
$ws.Range("M141").Value = -9056.155000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 165.81818
$ws.Range("I4").Value = 116.25
$ws.Range("J4").Value = 298
$ws.Range("K4").Value = 116.25
$ws.Range("L4").Value = 298
$ws.Range("M4").Value = -0.25
$ws.Range("N4").Value = -530
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()
$ws.Range("H10").Value = 13999.667
$ws.Range("J10").Value = 13999.667
$ws.Range("L10").Value = 13999.667
$ws.Range("N10").Value = -14339.667
$ws.Range("H37").Value = 2000
$ws.Range("I37").Value = 2000
$ws.Range("K37").Value = 2000
$ws.Range("M37").Value = -1727
$ws.Range("H44").Value = 39999.5
$ws.Range("J44").Value = 39999.5
$ws.Range("L44").Value = 39999.5
$ws.Range("N44").Value = -40975.5
$ws.Range("H55").Value = 30483
$ws.Range("I55").Value = 1450
$ws.Range("J55").Value = 44999.5
$ws.Range("K55").Value = 1450
$ws.Range("L55").Value = 44999.5
$ws.Range("M55").Value = -1135
$ws.Range("N55").Value = -45629.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2011.909
$ws.Range("I22").Value = 2183.1
$ws.Range("K22").Value = 2183.1
$ws.Range("M22").Value = -2010.1
$ws.Range("H35").Value = 69999
$ws.Range("J35").Value = 69999
$ws.Range("L35").Value = 69999
$ws.Range("N35").Value = -70619
$ws.Range("H80").Value = 1404.5834
$ws.Range("I80").Value = 1085.2
$ws.Range("J80").Value = 1632.7142
$ws.Range("K80").Value = 1085.2
$ws.Range("L80").Value = 1632.7142
$ws.Range("M80").Value = -87.20000000000005
$ws.Range("N80").Value = -3628.7142
$ws.Range("H83").Value = 1404.5834
$ws.Range("I83").Value = 1085.2
$ws.Range("J83").Value = 1632.7142
$ws.Range("K83").Value = 5426
$ws.Range("L83").Value = 8163.571
$ws.Range("M83").Value = -434
$ws.Range("N83").Value = -18147.571
$ws.Range("H88").Value = 14785.571
$ws.Range("J88").Value = 16550
$ws.Range("L88").Value = 16550
$ws.Range("N88").Value = -17362
$ws.Range("H91").Value = 14785.571
$ws.Range("J91").Value = 16550
$ws.Range("L91").Value = 16550
$ws.Range("N91").Value = -19358
$ws.Range("H134").Value = 6941.3213
$ws.Range("I134").Value = 6982.0415
$ws.Range("K134").Value = 20946.1245
$ws.Range("M134").Value = -18411.1245

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 15997.6
$ws.Range("H22").Value = 881.5714
$ws.Range("I22").Value = 593.6667
$ws.Range("J22").Value = 1399.8
$ws.Range("K22").Value = 593.6667
$ws.Range("L22").Value = 1399.8
$ws.Range("M22").Value = -243.6667
$ws.Range("N22").Value = -2099.8
$ws.Range("H31").Value = 3213.5334
$ws.Range("J31").Value = 3591.4
$ws.Range("L31").Value = 3591.4
$ws.Range("N31").Value = -4181.4
$ws.Range("H34").Value = 3213.5334
$ws.Range("J34").Value = 3591.4
$ws.Range("L34").Value = 3591.4
$ws.Range("N34").Value = -3995.4
$ws.Range("H105").Value = 2817.4167
$ws.Range("I105").Value = 1878.8
$ws.Range("K105").Value = 1878.8
$ws.Range("M105").Value = -131.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 656.8
$ws.Range("J23").Value = 861.3333
$ws.Range("L23").Value = 2583.9999
$ws.Range("N23").Value = -3053.9999
$ws.Range("H36").Value = 2085.4285
$ws.Range("J36").Value = 2266.6667
$ws.Range("L36").Value = 6800.000100000001
$ws.Range("N36").Value = -7138.000100000001
$ws.Range("H86").Value = 3516.5293
$ws.Range("I86").Value = 2746.875
$ws.Range("K86").Value = 8240.625
$ws.Range("M86").Value = -7054.625
$ws.Range("H89").Value = 3516.5293
$ws.Range("I89").Value = 2746.875
$ws.Range("K89").Value = 24721.875
$ws.Range("M89").Value = -18793.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 149.5
$ws.Range("I2").Value = 55
$ws.Range("K2").Value = 55
$ws.Range("M2").Value = 58
$ws.Range("H80").Value = 1357.2222
$ws.Range("I80").Value = 1249.25
$ws.Range("J80").Value = 1443.6
$ws.Range("K80").Value = 1249.25
$ws.Range("L80").Value = 1443.6
$ws.Range("M80").Value = -251.25
$ws.Range("N80").Value = -3439.6
$ws.Range("H83").Value = 1357.2222
$ws.Range("I83").Value = 1249.25
$ws.Range("J83").Value = 1443.6
$ws.Range("K83").Value = 6246.25
$ws.Range("L83").Value = 7218
$ws.Range("M83").Value = -1254.25
$ws.Range("N83").Value = -17202
$ws.Range("H122").Value = 1953.7368
$ws.Range("I122").Value = 1924.7646
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 5774.293799999999
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -3324.293799999999
$ws.Range("N122").Value = -11500

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 479.625
$ws.Range("J55").Value = 953.2
$ws.Range("L55").Value = 953.2
$ws.Range("N55").Value = -1299.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 69476
$ws.Range("J99").Value = 69476
$ws.Range("L99").Value = 69476
$ws.Range("N99").Value = -75466
